$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2498276264660149
$ws.Cells.Item(2, 3).Value = 0.03160055733154365
$ws.Cells.Item(2, 5).Value = 0.495341007136858
$ws.Cells.Item(2, 6).Value = 2.177157614419173
$ws.Cells.Item(2, 7).Value = 0.002439468186846574
$ws.Cells.Item(2, 9).Value = 0.5704123665821079
$ws.Cells.Item(2, 10).Value = 0.0495201053083596
$ws.Cells.Item(2, 11).Value = 0.2501965568839637
$ws.Cells.Item(2, 13).Value = 0.4114831580304923
$ws.Cells.Item(2, 14).Value = 1.553447452456883
$ws.Cells.Item(2, 15).Value = 2.588490551675477

$ws.Cells.Item(3, 2).Value = 0.2211950567261738
$ws.Cells.Item(3, 3).Value = 0.02770956227770682
$ws.Cells.Item(3, 5).Value = 0.4859447645738868
$ws.Cells.Item(3, 6).Value = 2.165079812083405
$ws.Cells.Item(3, 7).Value = 0.002441512221929623
$ws.Cells.Item(3, 9).Value = 0.5758158946234104
$ws.Cells.Item(3, 10).Value = 0.04977538768123857
$ws.Cells.Item(3, 11).Value = 0.2198446522830579
$ws.Cells.Item(3, 13).Value = 0.3912740604080369
$ws.Cells.Item(3, 14).Value = 1.569116659481365
$ws.Cells.Item(3, 15).Value = 2.606388713964037

$ws.Cells.Item(4, 2).Value = 0.2036119042027735
$ws.Cells.Item(4, 3).Value = 0.02530892200260837
$ws.Cells.Item(4, 5).Value = 0.4804096852856716
$ws.Cells.Item(4, 6).Value = 2.158814390840405
$ws.Cells.Item(4, 7).Value = 0.002442834231795221
$ws.Cells.Item(4, 9).Value = 0.5794128013568454
$ws.Cells.Item(4, 10).Value = 0.04994204892958987
$ws.Cells.Item(4, 11).Value = 0.2011851017825279
$ws.Cells.Item(4, 13).Value = 0.3790354067799555
$ws.Cells.Item(4, 14).Value = 1.579241844993279
$ws.Cells.Item(4, 15).Value = 2.618731963483526

$ws.Cells.Item(5, 2).Value = 0.196446407544812
$ws.Cells.Item(5, 3).Value = 0.02432777136468189
$ws.Cells.Item(5, 5).Value = 0.4782130792092616
$ws.Cells.Item(5, 6).Value = 2.156550586709713
$ws.Cells.Item(5, 7).Value = 0.002443389849031538
$ws.Cells.Item(5, 9).Value = 0.5809487545506364
$ws.Cells.Item(5, 10).Value = 0.05001246310872443
$ws.Cells.Item(5, 11).Value = 0.1935757378393674
$ws.Cells.Item(5, 13).Value = 0.3740909829490207
$ws.Cells.Item(5, 14).Value = 1.583494755514323
$ws.Cells.Item(5, 15).Value = 2.624102330813784

$ws.Cells.Item(6, 2).Value = 0.1952565821264614
$ws.Cells.Item(6, 3).Value = 0.02416467972764735
$ws.Cells.Item(6, 5).Value = 0.4778518992579635
$ws.Cells.Item(6, 6).Value = 2.156192167008825
$ws.Cells.Item(6, 7).Value = 0.002443483130391438
$ws.Cells.Item(6, 9).Value = 0.5812080373221775
$ws.Cells.Item(6, 10).Value = 0.05002430636270194
$ws.Cells.Item(6, 11).Value = 0.1923118924588749
$ws.Cells.Item(6, 13).Value = 0.3732725644710229
$ws.Cells.Item(6, 14).Value = 1.584208608103149
$ws.Cells.Item(6, 15).Value = 2.62501463493733

$ws.Cells.Item(7, 2).Value = 0.2035152681071679
$ws.Cells.Item(7, 3).Value = 0.02529570143748572
$ws.Cells.Item(7, 5).Value = 0.4803798221604367
$ws.Cells.Item(7, 6).Value = 2.158782688423727
$ws.Cells.Item(7, 7).Value = 0.002442841656862926
$ws.Cells.Item(7, 9).Value = 0.5794332315875437
$ws.Cells.Item(7, 10).Value = 0.04994298843783085
$ws.Cells.Item(7, 11).Value = 0.2010825007208012
$ws.Cells.Item(7, 13).Value = 0.3789685504425222
$ws.Cells.Item(7, 14).Value = 1.579298687731674
$ws.Cells.Item(7, 15).Value = 2.618803011908696

$ws.Cells.Item(8, 2).Value = 0.2399559218937952
$ws.Cells.Item(8, 3).Value = 0.03026136004650937
$ws.Cells.Item(8, 5).Value = 0.492052607244247
$ws.Cells.Item(8, 6).Value = 2.172754490664531
$ws.Cells.Item(8, 7).Value = 0.002440159101442501
$ws.Cells.Item(8, 9).Value = 0.5722175707742032
$ws.Cells.Item(8, 10).Value = 0.04960607126662886
$ws.Cells.Item(8, 11).Value = 0.2397363157954118
$ws.Cells.Item(8, 13).Value = 0.4044799448189949
$ws.Cells.Item(8, 14).Value = 1.558745557199124
$ws.Cells.Item(8, 15).Value = 2.59438087013605

$ws.Cells.Item(9, 2).Value = 0.3113789172515453
$ws.Cells.Item(9, 3).Value = 0.03990624049919234
$ws.Cells.Item(9, 5).Value = 0.5167998296414353
$ws.Cells.Item(9, 6).Value = 2.209279702729887
$ws.Cells.Item(9, 7).Value = 0.00243542769968483
$ws.Cells.Item(9, 9).Value = 0.5602824866228353
$ws.Cells.Item(9, 10).Value = 0.04902384916616853
$ws.Cells.Item(9, 11).Value = 0.3153357970497837
$ws.Cells.Item(9, 13).Value = 0.4558480420676432
$ws.Cells.Item(9, 14).Value = 1.5224401873856
$ws.Cells.Item(9, 15).Value = 2.557231651127097

$ws.Cells.Item(10, 2).Value = 0.3638142937143698
$ws.Cells.Item(10, 3).Value = 0.04693498526765438
$ws.Cells.Item(10, 5).Value = 0.5361138400756289
$ws.Cells.Item(10, 6).Value = 2.241683164328151
$ws.Cells.Item(10, 7).Value = 0.00243227092562702
$ws.Cells.Item(10, 9).Value = 0.5528641946998292
$ws.Cells.Item(10, 10).Value = 0.0486436292249941
$ws.Cells.Item(10, 11).Value = 0.3707413869959737
$ws.Cells.Item(10, 13).Value = 0.4943999989905663
$ws.Cells.Item(10, 14).Value = 1.498201717990623
$ws.Cells.Item(10, 15).Value = 2.536490769865452

$ws.Cells.Item(11, 2).Value = 0.3876568086297141
$ws.Cells.Item(11, 3).Value = 0.05011996108888184
$ws.Cells.Item(11, 5).Value = 0.5451463390763536
$ws.Cells.Item(11, 6).Value = 2.257634870087529
$ws.Cells.Item(11, 7).Value = 0.002430903508278237
$ws.Cells.Item(11, 9).Value = 0.5497827126273691
$ws.Cells.Item(11, 10).Value = 0.0484809161665094
$ws.Cells.Item(11, 11).Value = 0.3959141129836894
$ws.Cells.Item(11, 13).Value = 0.5121136901543295
$ws.Cells.Item(11, 14).Value = 1.487703147732876
$ws.Cells.Item(11, 15).Value = 2.528479009195252

$ws.Cells.Item(12, 2).Value = 0.3966834376075212
$ws.Cells.Item(12, 3).Value = 0.05132421382975849
$ws.Cells.Item(12, 5).Value = 0.5486021124177967
$ws.Cells.Item(12, 6).Value = 2.263849523063953
$ws.Cells.Item(12, 7).Value = 0.002430395518921069
$ws.Cells.Item(12, 9).Value = 0.5486579979159494
$ws.Cells.Item(12, 10).Value = 0.04842077024235891
$ws.Cells.Item(12, 11).Value = 0.4054414652866853
$ws.Cells.Item(12, 13).Value = 0.5188465726539278
$ws.Cells.Item(12, 14).Value = 1.483803471174067
$ws.Cells.Item(12, 15).Value = 2.525649908286027

$ws.Cells.Item(13, 2).Value = 0.3947394889169686
$ws.Cells.Item(13, 3).Value = 0.05106493841991266
$ws.Cells.Item(13, 5).Value = 0.547856278314157
$ws.Cells.Item(13, 6).Value = 2.262503345044365
$ws.Cells.Item(13, 7).Value = 0.002430504487584946
$ws.Cells.Item(13, 9).Value = 0.5488983492099209
$ws.Cells.Item(13, 10).Value = 0.04843365843020386
$ws.Cells.Item(13, 11).Value = 0.4033898100094575
$ws.Cells.Item(13, 13).Value = 0.5173954139039694
$ws.Cells.Item(13, 14).Value = 1.48463995990865
$ws.Cells.Item(13, 15).Value = 2.526250096590019

$ws.Cells.Item(14, 2).Value = 0.3883994779955913
$ws.Cells.Item(14, 3).Value = 0.05021907256239899
$ws.Cells.Item(14, 5).Value = 0.5454299393720703
$ws.Cells.Item(14, 6).Value = 2.258142664813633
$ws.Cells.Item(14, 7).Value = 0.002430861519070922
$ws.Cells.Item(14, 9).Value = 0.5496893362753354
$ws.Cells.Item(14, 10).Value = 0.04847593848832155
$ws.Cells.Item(14, 11).Value = 0.3966980371754971
$ws.Cells.Item(14, 13).Value = 0.5126671069842388
$ws.Cells.Item(14, 14).Value = 1.487380797573815
$ws.Cells.Item(14, 15).Value = 2.528242152711101

$ws.Cells.Item(15, 2).Value = 0.3845157614343293
$ws.Cells.Item(15, 3).Value = 0.04970071598521031
$ws.Cells.Item(15, 5).Value = 0.5439483395545466
$ws.Cells.Item(15, 6).Value = 2.255494291236474
$ws.Cells.Item(15, 7).Value = 0.002431081489100648
$ws.Cells.Item(15, 9).Value = 0.5501793318962918
$ws.Cells.Item(15, 10).Value = 0.04850202755793198
$ws.Cells.Item(15, 11).Value = 0.3925984661887014
$ws.Cells.Item(15, 13).Value = 0.5097741428565072
$ws.Cells.Item(15, 14).Value = 1.489069525723309
$ws.Cells.Item(15, 15).Value = 2.52948901584594

$ws.Cells.Item(16, 2).Value = 0.3622558741229795
$ws.Cells.Item(16, 3).Value = 0.04672658593061385
$ws.Cells.Item(16, 5).Value = 0.5355284977012076
$ws.Cells.Item(16, 6).Value = 2.240665048447767
$ws.Cells.Item(16, 7).Value = 0.002432361666414394
$ws.Cells.Item(16, 9).Value = 0.5530714774161254
$ws.Cells.Item(16, 10).Value = 0.04865446883003699
$ws.Cells.Item(16, 11).Value = 0.3690956145704831
$ws.Cells.Item(16, 13).Value = 0.4932458923571303
$ws.Cells.Item(16, 14).Value = 1.498898444159622
$ws.Cells.Item(16, 15).Value = 2.537043007664892

$ws.Cells.Item(17, 2).Value = 0.3485970907538558
$ws.Cells.Item(17, 3).Value = 0.04489883979380238
$ws.Cells.Item(17, 5).Value = 0.5304262656850511
$ws.Cells.Item(17, 6).Value = 2.23187794963853
$ws.Cells.Item(17, 7).Value = 0.002433164555925076
$ws.Cells.Item(17, 9).Value = 0.5549208105705041
$ws.Cells.Item(17, 10).Value = 0.04875060919838248
$ws.Cells.Item(17, 11).Value = 0.3546689563291636
$ws.Cells.Item(17, 13).Value = 0.4831513027131251
$ws.Cells.Item(17, 14).Value = 1.505063327756815
$ws.Cells.Item(17, 15).Value = 2.542041767372254

$ws.Cells.Item(18, 2).Value = 0.3407399512793745
$ws.Cells.Item(18, 3).Value = 0.04384640048323263
$ws.Cells.Item(18, 5).Value = 0.5275147993502856
$ws.Cells.Item(18, 6).Value = 2.226937852430524
$ws.Cells.Item(18, 7).Value = 0.002433632817698999
$ws.Cells.Item(18, 9).Value = 0.5560120880587895
$ws.Cells.Item(18, 10).Value = 0.04880687170175158
$ws.Cells.Item(18, 11).Value = 0.346368182396958
$ws.Cells.Item(18, 13).Value = 0.4773617625156561
$ws.Cells.Item(18, 14).Value = 1.508658866821294
$ws.Cells.Item(18, 15).Value = 2.545050890079096

$ws.Cells.Item(19, 2).Value = 0.3380795074293985
$ws.Cells.Item(19, 3).Value = 0.04348986322600012
$ws.Cells.Item(19, 5).Value = 0.5265330143273275
$ws.Cells.Item(19, 6).Value = 2.225284804157525
$ws.Cells.Item(19, 7).Value = 0.002433792474035123
$ws.Cells.Item(19, 9).Value = 0.5563863137704921
$ws.Cells.Item(19, 10).Value = 0.04882608710123826
$ws.Cells.Item(19, 11).Value = 0.3435571922130976
$ws.Cells.Item(19, 13).Value = 0.4754043823994607
$ws.Cells.Item(19, 14).Value = 1.509884782994206
$ws.Cells.Item(19, 15).Value = 2.546092733474524

$ws.Cells.Item(20, 2).Value = 0.350051195348243
$ws.Cells.Item(20, 3).Value = 0.04509352767341568
$ws.Cells.Item(20, 5).Value = 0.5309670061819673
$ws.Cells.Item(20, 6).Value = 2.232801552330869
$ws.Cells.Item(20, 7).Value = 0.002433078418622383
$ws.Cells.Item(20, 9).Value = 0.5547210902735387
$ws.Cells.Item(20, 10).Value = 0.04874027503710643
$ws.Cells.Item(20, 11).Value = 0.35620500662219
$ws.Cells.Item(20, 13).Value = 0.4842241726836889
$ws.Cells.Item(20, 14).Value = 1.504401925417563
$ws.Cells.Item(20, 15).Value = 2.541495774920321

$ws.Cells.Item(21, 2).Value = 0.3902617514470705
$ws.Cells.Item(21, 3).Value = 0.05046757381421685
$ws.Cells.Item(21, 5).Value = 0.5461416548465792
$ws.Cells.Item(21, 6).Value = 2.259418778255167
$ws.Cells.Item(21, 7).Value = 0.002430756383983227
$ws.Cells.Item(21, 9).Value = 0.5494558593803625
$ws.Cells.Item(21, 10).Value = 0.04846347994793732
$ws.Cells.Item(21, 11).Value = 0.3986637128982693
$ws.Cells.Item(21, 13).Value = 0.5140552463224424
$ws.Cells.Item(21, 14).Value = 1.486573686831891
$ws.Cells.Item(21, 15).Value = 2.527651479099063

$ws.Cells.Item(22, 2).Value = 0.4165297082394375
$ws.Cells.Item(22, 3).Value = 0.05396914717950096
$ws.Cells.Item(22, 5).Value = 0.5562652452059922
$ws.Cells.Item(22, 6).Value = 2.277829354873134
$ws.Cells.Item(22, 7).Value = 0.002429296030628084
$ws.Cells.Item(22, 9).Value = 0.5462605695099967
$ws.Cells.Item(22, 10).Value = 0.0482911446405847
$ws.Cells.Item(22, 11).Value = 0.42638349348681
$ws.Cells.Item(22, 13).Value = 0.5336977182378959
$ws.Cells.Item(22, 14).Value = 1.475364346244845
$ws.Cells.Item(22, 15).Value = 2.519797086657292

$ws.Cells.Item(23, 2).Value = 0.4025112619822551
$ws.Cells.Item(23, 3).Value = 0.0521012817728348
$ws.Cells.Item(23, 5).Value = 0.550843264131359
$ws.Cells.Item(23, 6).Value = 2.267910467559531
$ws.Cells.Item(23, 7).Value = 0.002430070226750367
$ws.Cells.Item(23, 9).Value = 0.5479434535562646
$ws.Cells.Item(23, 10).Value = 0.04838234076625803
$ws.Cells.Item(23, 11).Value = 0.4115917755200087
$ws.Cells.Item(23, 13).Value = 0.5232008712082532
$ws.Cells.Item(23, 14).Value = 1.48130648754805
$ws.Cells.Item(23, 15).Value = 2.523879871079345

$ws.Cells.Item(24, 2).Value = 0.3493938087122785
$ws.Cells.Item(24, 3).Value = 0.04500551440095535
$ws.Cells.Item(24, 5).Value = 0.5307224692593451
$ws.Cells.Item(24, 6).Value = 2.232383643536238
$ws.Cells.Item(24, 7).Value = 0.002433117340504283
$ws.Cells.Item(24, 9).Value = 0.5548112963810503
$ws.Cells.Item(24, 10).Value = 0.04874494402718454
$ws.Cells.Item(24, 11).Value = 0.3555105791061521
$ws.Cells.Item(24, 13).Value = 0.4837390846247729
$ws.Cells.Item(24, 14).Value = 1.504700785726747
$ws.Cells.Item(24, 15).Value = 2.541742196757596

$ws.Cells.Item(25, 2).Value = 0.2920627771524664
$ws.Cells.Item(25, 3).Value = 0.03730705825896052
$ws.Cells.Item(25, 5).Value = 0.5099062328966539
$ws.Cells.Item(25, 6).Value = 2.1984211292541
$ws.Cells.Item(25, 7).Value = 0.002436651354954738
$ws.Cells.Item(25, 9).Value = 0.5632741122278517
$ws.Cells.Item(25, 10).Value = 0.04917298583252894
$ws.Cells.Item(25, 11).Value = 0.3153357970497837
$ws.Cells.Item(25, 13).Value = 0.4418086634435383
$ws.Cells.Item(25, 14).Value = 1.53183392156012
$ws.Cells.Item(25, 15).Value = 2.566131056362266

Write-Host "Updated 264 cells"